$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.411.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "'1.995.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.77%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'330.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.38%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4934"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.05%  "
$ws.Range("D8").Value = "'0.4176"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.28%  "
$ws.Range("D9").Value = "'53.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'0.08799"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.06%  "
$ws.Range("D11").Value = "'1.111"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.13%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'23.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.47%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'2.061.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'8.018"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.57%  "
$ws.Range("D15").Value = "'6.460"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.61%  "
$ws.Range("D16").Value = "'96.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.48%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'0.00001105"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").Value = "'0.06623"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "'19.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.28%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'5.957"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("D23").Value = "'29.470.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "'11.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.28%  "
$ws.Range("D25").Value = "'2.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value = "'2.342.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'6.676"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "'157.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").Value = "'20.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.47%  "
$ws.Range("D30").Value = "'2.345"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.55%  "
$ws.Range("D31").Value = "'126.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.83%  "
$ws.Range("D33").Value = "'0.09904"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -13.40%  "
$ws.Range("D35").Value = "'5.815"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.07%  "
$ws.Range("D36").Value = "'3.779"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.74%  "
$ws.Range("D37").Value = "'9.570"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.45%  "
$ws.Range("D38").Value = "'0.02443"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.06358"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.26%  "
$ws.Range("D40").Value = "'1.278"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6489"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.42%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.01%  "
$ws.Range("D43").Value = "'0.2062"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.02%  "
$ws.Range("D44").Value = "'1.006"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.6317"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.71%  "
$ws.Range("D46").Value = "'13.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.55%  "
$ws.Range("D47").Value = "'2.193"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.51%  "
$ws.Range("D48").Value = "'1.262"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("D49").Value = "'3.540"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000329"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.86%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.06982"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.94%  "
